$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.664.85'
$ws.Range("E2").Value = '  +4.55%  '

$ws.Range("D3").Value = '3.496.15'
$ws.Range("E3").Value = '  +2.68%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '590.82'
$ws.Range("E5").Value = '  +3.74%  '

$ws.Range("D6").Value = '168.89'
$ws.Range("E6").Value = '  +4.39%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  +8.95%  '

$ws.Range("D9").Value = '3.493.79'
$ws.Range("E9").Value = '  +2.46%  '

$ws.Range("D10").Value = '0.129'
$ws.Range("E10").Value = '  +7.81%  '

$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("E12").Value = '  +4.36%  '

$ws.Range("D13").Value = '4.102.36'
$ws.Range("E13").Value = '  +2.40%  '

$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("D15").Value = '28.23'
$ws.Range("E15").Value = '  +5.08%  '

$ws.Range("E16").Value = '  +4.28%  '

$ws.Range("D17").Value = '66.653.61'
$ws.Range("E17").Value = '  +4.33%  '

$ws.Range("D18").Value = '3.502.17'
$ws.Range("E18").Value = '  +1.85%  '

$ws.Range("D19").Value = '6.33'
$ws.Range("E19").Value = '  +3.87%  '

$ws.Range("D20").Value = '14.18'
$ws.Range("E20").Value = '  +4.57%  '

$ws.Range("D21").Value = '392.74'
$ws.Range("E21").Value = '  +4.55%  '

$ws.Range("D22").Value = '7.97'
$ws.Range("E22").Value = '  +3.04%  '

$ws.Range("D23").Value = '73.46'
$ws.Range("E23").Value = '  +3.55%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("E25").Value = '  +4.63%  '

$ws.Range("E26").Value = '  +6.26%  '

$ws.Range("D27").Value = '10.19'
$ws.Range("E27").Value = '  +7.33%  '

$ws.Range("E28").Value = '  +2.59%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = '6.36'
$ws.Range("E30").Value = '  +4.88%  '

$ws.Range("E31").Value = '  +6.73%  '

$ws.Range("E32").Value = '  +3.66%  '

$ws.Range("D33").Value = '23.63'
$ws.Range("E33").Value = '  +3.54%  '

$ws.Range("E34").Value = '  +5.22%  '

$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  +9.90%  '

$ws.Range("D37").Value = '162.68'
$ws.Range("E37").Value = '  +1.79%  '

$ws.Range("D38").Value = '0.884'
$ws.Range("E38").Value = '  +2.91%  '

$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  +6.94%  '

$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Value = '4.69'
$ws.Range("E40").Value = '  +7.19%  '

$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '27.56'
$ws.Range("E41").Value = '  +6.21%  '

$ws.Range("E42").Value = '  +3.09%  '

$ws.Range("D43").Value = '6.75'
$ws.Range("E43").Value = '  +5.34%  '

$ws.Range("D44").Value = '26.50'
$ws.Range("E44").Value = '  +2.73%  '

$ws.Range("D45").Value = '2.788.72'
$ws.Range("E45").Value = '  +0.66%  '

$ws.Range("D46").Value = '43.19'
$ws.Range("E46").Value = '  +1.13%  '

$ws.Range("D47").Value = '0.0314'
$ws.Range("E47").Value = '  +2.80%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '352.45'
$ws.Range("E48").Value = '  +7.42%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '2.50'
$ws.Range("E49").Value = '  +4.18%  '

$ws.Range("E50").Value = '  +6.04%  '

$ws.Range("D51").Value = '33.93'
$ws.Range("E51").Value = '  +13.79%  '
